$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '28.491.53'
    'E2' = '  -0.10%  '
    'D3' = '1.825.11'
    'E3' = '  -0.30%  '
    'E4' = '  +0.32%  '
    'D5' = '316.50'
    'E5' = '  +0.31%  '
    'E6' = '  +0.23%  '
    'D7' = '0.5178'
    'E7' = '  +2.07%  '
    'D8' = '0.3863'
    'E8' = '  -1.55%  '
    'D9' = '0.08408'
    'E9' = '  +8.85%  '
    'D10' = '1.119'
    'E10' = '  +0.50%  '
    'D11' = '41.96'
    'E11' = '  +0.07%  '
    'D12' = '6.414'
    'E12' = '  +2.42%  '
    'D13' = '21.20'
    'E13' = '  +0.66%  '
    'D14' = '1.004'
    'E14' = '  +0.17%  '
    'D15' = '7.504'
    'E15' = '  -0.56%  '
    'D16' = '1.819.22'
    'E16' = '  -0.05%  '
    'D17' = '94.24'
    'E17' = '  +0.76%  '
    'D18' = '0.00001130'
    'E18' = '  +4.33%  '
    'D19' = '0.06638'
    'D20' = '17.75'
    'E20' = '  -0.06%  '
    'E21' = '  +0.22%  '
    'D22' = '6.077'
    'E22' = '  -0.92%  '
    'D23' = '28.540.55'
    'E23' = '  +0.02%  '
    'D24' = '11.41'
    'E24' = '  +2.40%  '
    'D25' = '2.291'
    'E25' = '  +1.62%  '
    'D26' = '21.22'
    'E26' = '  +2.77%  '
    'D27' = '159.52'
    'E27' = '  +1.64%  '
    'D28' = '2.033.18'
    'E28' = '  -0.13%  '
    'D29' = '2.400'
    'E29' = '  -1.26%  '
    'D30' = '125.94'
    'E30' = '  +0.59%  '
    'D31' = '0.1095'
    'E31' = '  +0.35%  '
    'E32' = '  -2.99%  '
    'D33' = '0.07686'
    'E33' = '  +8.14%  '
    'D34' = '5.741'
    'E34' = '  +1.36%  '
    'D35' = '3.675'
    'E35' = '  +0.22%  '
    'E36' = '  +0.64%  '
    'D37' = '0.02382'
    'E37' = '  +2.53%  '
    'D38' = '5.289'
    'E38' = '  +3.12%  '
    'D39' = '8.758'
    'E39' = '  -2.51%  '
    'D40' = '0.6417'
    'E40' = '  +2.78%  '
    'D41' = '11.53'
    'E41' = '  +2.78%  '
    'D42' = '1.196'
    'E42' = '  +0.46%  '
    'D43' = '1.400'
    'E43' = '  +0.16%  '
    'D44' = '13.58'
    'E44' = '  +0.84%  '
    'D45' = '0.6149'
    'E45' = '  +4.20%  '
    'D46' = '3.796'
    'E46' = '  +2.04%  '
    'D47' = '127.98'
    'E47' = '  +2.89%  '
    'D48' = '1.999'
    'E48' = '  +1.38%  '
    'D49' = '1.206'
    'E49' = '  +1.91%  '
    'D50' = '0.06996'
    'E50' = '  +0.98%  '
    'D51' = '74.42'
    'E51' = '  +0.81%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates"